# Fix 5bus esd1 case
# Adds a 5th generator "PV_2" (on Bus 1) to the PJM 5-bus ESD1 case, wires it
# up through GCost / SFRCost / SRCost / NSRCost, re-points the ESD1_1 storage
# unit at the new generator/bus, and refreshes the EDTSlot "ug" mask to cover
# 5 units instead of 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# PV sheet: add PV_2 (uid=3) on bus 1
# ---------------------------------------------------------------------------
$wsPV = $wb.Worksheets.Item("PV")
$wsPV.Range("A4:U4").Copy()
$wsPV.Range("A5:U5").PasteSpecial(-4122)

$wsPV.Range("A5").Value = 3
$wsPV.Range("B5").Value = "PV_2"
$wsPV.Range("C5").Value = 1
$wsPV.Range("D5").Value = "PV 2"
$wsPV.Range("E5").Value = 100
$wsPV.Range("F5").Value = 230
$wsPV.Range("G5").Value = 1
$wsPV.Range("I5").Value = 0
$wsPV.Range("J5").Value = 0
$wsPV.Range("K5").Value = 99
$wsPV.Range("L5").Value = -99
$wsPV.Range("M5").Value = 99
$wsPV.Range("N5").Value = -99
$wsPV.Range("O5").Value = 1
$wsPV.Range("P5").Value = 1.4
$wsPV.Range("Q5").Value = 0.6
$wsPV.Range("R5").Value = 0.5
$wsPV.Range("S5").Value = 0
$wsPV.Range("T5").Value = 0.01
$wsPV.Range("U5").Value = 0.3

# ---------------------------------------------------------------------------
# ESD1 sheet: point ESD1_1 at bus 1 / gen PV_2 (was bus 0 / gen PV_1)
# ---------------------------------------------------------------------------
$wsESD1 = $wb.Worksheets.Item("ESD1")
$wsESD1.Range("E2").Value = 1
$wsESD1.Range("F2").Value = "PV_2"

# ---------------------------------------------------------------------------
# GCost sheet: add GCost_5 row for PV_2
# ---------------------------------------------------------------------------
$wsGCost = $wb.Worksheets.Item("GCost")
$wsGCost.Range("A5:K5").Copy()
$wsGCost.Range("A6:K6").PasteSpecial(-4122)

$wsGCost.Range("A6").Value = 4
$wsGCost.Range("B6").Value = "GCost_5"
$wsGCost.Range("C6").Value = 1
$wsGCost.Range("D6").Value = "GCost_5"
$wsGCost.Range("E6").Value = "PV_2"
$wsGCost.Range("F6").Value = 2
$wsGCost.Range("G6").Value = 0
$wsGCost.Range("H6").Value = 0
$wsGCost.Range("I6").Value = 0
$wsGCost.Range("J6").Value = 0
$wsGCost.Range("K6").Value = 0

# ---------------------------------------------------------------------------
# SFRCost sheet: add SFRC_5 row for PV_2
# ---------------------------------------------------------------------------
$wsSFRCost = $wb.Worksheets.Item("SFRCost")
$wsSFRCost.Range("A5:E5").Copy()
$wsSFRCost.Range("A6:E6").PasteSpecial(-4122)

$wsSFRCost.Range("A6").Value = "4"
$wsSFRCost.Range("B6").Value = "SFRC_5"
$wsSFRCost.Range("C6").Value = "PV_2"
$wsSFRCost.Range("D6").Value = "0"
$wsSFRCost.Range("E6").Value = "0"

# ---------------------------------------------------------------------------
# SRCost sheet: add SRC_5 row for PV_2
# ---------------------------------------------------------------------------
$wsSRCost = $wb.Worksheets.Item("SRCost")
$wsSRCost.Range("A5:D5").Copy()
$wsSRCost.Range("A6:D6").PasteSpecial(-4122)

$wsSRCost.Range("A6").Value = "4"
$wsSRCost.Range("B6").Value = "SRC_5"
$wsSRCost.Range("C6").Value = "PV_2"
$wsSRCost.Range("D6").Value = 0.1

# ---------------------------------------------------------------------------
# NSRCost sheet: add NSRC_5 row for PV_2
# ---------------------------------------------------------------------------
$wsNSRCost = $wb.Worksheets.Item("NSRCost")
$wsNSRCost.Range("A5:D5").Copy()
$wsNSRCost.Range("A6:D6").PasteSpecial(-4122)

$wsNSRCost.Range("A6").Value = "4"
$wsNSRCost.Range("B6").Value = "NSRC_5"
$wsNSRCost.Range("C6").Value = "PV_2"
$wsNSRCost.Range("D6").Value = 0.1

# ---------------------------------------------------------------------------
# EDTSlot sheet: "ug" mask now covers 5 units instead of 4
# ---------------------------------------------------------------------------
$wsEDTSlot = $wb.Worksheets.Item("EDTSlot")
$wsEDTSlot.Range("E2").Value = "1,1,1,1,1"
$wsEDTSlot.Range("E3").Value = "1,1,1,1,1"
$wsEDTSlot.Range("E4").Value = "1,1,1,1,1"
$wsEDTSlot.Range("E5").Value = "1,1,1,1,1"

# ---------------------------------------------------------------------------
# View-state touch-ups (selections / zoom / active sheet) to mirror the
# author's final on-screen state.
# ---------------------------------------------------------------------------
$wsEDTSlot.Select()
$wsEDTSlot.Range("I12").Select()

$wsSFRCost.Select()
$excel.ActiveWindow.Zoom = 207
$wsSFRCost.Range("C7").Select()

$wsSRCost.Select()
$excel.ActiveWindow.Zoom = 230
$wsSRCost.Range("C7").Select()

$wsNSRCost.Select()
$excel.ActiveWindow.Zoom = 200
$wsNSRCost.Range("C7").Select()

$wsPV.Select()
$excel.ActiveWindow.Zoom = 140
$wsPV.Range("K5").Select()

$wsGCost.Select()
$wsGCost.Range("F29").Select()

$wsLine = $wb.Worksheets.Item("Line")
$wsLine.Select()
$wsLine.Range("N20").Select()

$wsESD1.Select()
$wsESD1.Range("F3").Select()
